# "update core tms 1&2"
# The case id for the first (core) test-management row on Sheet1 is
# refreshed to a newly generated value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "CA-HGTFBRXD"
